# Add the "last_coupon_date" worksheet as the last (3rd) tab and populate
# it with bond names (copied from the close_price header row) and their
# last-coupon dates, formatted to match the workbook's existing look.

$wb = $excel.ActiveWorkbook

# --- add the new sheet after the last existing sheet -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "last_coupon_date"

# --- grab the existing "Times New Roman" style used on close_price ----
$closePrice = $wb.Worksheets.Item("close_price")
$styleSrc = $closePrice.Range("A1")

# --- header row --------------------------------------------------------
$ws.Range("A1").Value = "Bond"
$ws.Range("B1").Value = "Last coupon date"

$styleSrc.Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A1:B1").Font.Size = 12
$ws.Range("B1").HorizontalAlignment = -4152   # xlRight

# --- bond names (column A) + last coupon dates (column B) -------------
$bonds = @(
    "CAN 0.5 Mar 1 2022",
    "CAN 2.75 Jun 1 2022",
    "CAN 1.75 Mar 1 2023",
    "CAN 1.5 Jun 1 2023",
    "CAN 2.25 Mar 1 2024",
    "CAN 2.5 Jun 1 2024",
    "CAN 1.25 Mar 1 2025",
    "CAN 2.25 Jun 1 2025",
    "CAN 0.25 Mar 1 2026",
    "CAN 1.5 Jun 1 2026"
)
$dates = @(44440, 44531, 44440, 44531, 44440, 44531, 44440, 44531, 44440, 44531)

for ($i = 0; $i -lt $bonds.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $bonds[$i]
    $ws.Range("B$row").Value = $dates[$i]
    $ws.Rows.Item($row).RowHeight = 30
}

$styleSrc.Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)   # xlPasteFormats (Times New Roman, justify)

$ws.Range("B2:B11").NumberFormat = "yyyy-mm-dd"
$ws.Range("B2:B11").Font.Name = "Times New Roman"
$ws.Range("B2:B11").Font.Size = 12

$ws.Columns.Item(2).ColumnWidth = 16

# --- sheet view: zoom + selection --------------------------------------
$ws.Select()
$ws.Range("A2").Select()
$excel.ActiveWindow.Zoom = 223

# --- close_price: selection moves off the header -----------------------
$closePrice.Select()
$closePrice.Range("E6").Select()

# --- make the new sheet the active tab ---------------------------------
$ws.Select()
